# Generate Report for Archive
# Status text flips from "Ready for handoff" to "In Translation" on every
# sheet that reports it, and the "Status" column narrows to fit the
# shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet keeps one status column per locale (E = zh-cn, F = de-de)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-locale detail sheets keep the status in column C
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# The Status columns auto-fit narrower now that "In Translation" (14 chars)
# is shorter than "Ready for handoff" (17 chars).
$overview.Columns.Item(5).ColumnWidth = 12.45
$overview.Columns.Item(6).ColumnWidth = 12.45
$zhcn.Columns.Item(3).ColumnWidth = 12.45
$dede.Columns.Item(3).ColumnWidth = 12.45
